$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Report"
$ws.Name = "Report"

# Update the membership counts in column B for each society
$ws.Range("B2").Value = 4100
$ws.Range("B3").Value = 3200
$ws.Range("B4").Value = 1200
$ws.Range("B5").Value = 1760
$ws.Range("B6").Value = 4000

# Restore the selection/cursor position to E10, matching the saved view state
$ws.Range("E10").Select()
